{"js": "const pairs = [\n  [\"37-31=\", \"83-9=\"],\n  [\"42-22=\", \"80+16=\"],\n  [\"42+35=\", \"74-23=\"],\n  [\"99-55=\", \"91-61=\"],\n  [\"8+19=\", \"16+20=\"],\n  [\"70-32=\", \"11+3=\"],\n  [\"52+15=\", \"12+23=\"],\n  [\"16+46=\", \"28+1=\"],\n  [\"4+80=\", \"24+62=\"],\n  [\"52+18=\", \"31+58=\"],\n  [\"42+41=\", \"56-49=\"],\n  [\"11+86=\", \"71+23=\"],\n  [\"23-14=\", \"67-14=\"],\n  [\"69-33=\", \"53-25=\"],\n  [\"56+39=\", \"99-11=\"],\n  [\"27+13=\", \"20+34=\"],\n  [\"68+6=\", \"8+67=\"],\n  [\"74+19=\", \"84-83=\"],\n  [\"44-11=\", \"47+7=\"],\n  [\"62-21=\", \"73-45=\"],\n  [\"90-69=\", \"47-32=\"],\n  [\"6+60=\", \"81-78=\"],\n  [\"58-29=\", \"81-74=\"],\n  [\"76+11=\", \"4+9=\"],\n  [\"81-15=\", \"19+26=\"],\n  [\"90-6=\", \"13+53=\"],\n  [\"8+44=\", \"87-37=\"],\n  [\"42+12=\", \"29+4=\"],\n  [\"77+3=\", \"36-19=\"],\n  [\"77-35=\", \"58-54=\"],\n  [\"85-44=\", \"83-19=\"],\n  [\"66+28=\", \"89-2=\"],\n  [\"37-12=\", \"73-41=\"],\n  [\"3+3=\", \"52-50=\"],\n  [\"13+76=\", \"40+59=\"],\n  [\"56+23=\", \"54+2=\"],\n  [\"27+71=\", \"21+0=\"],\n  [\"53+25=\", \"8+51=\"],\n  [\"64-44=\", \"61-50=\"],\n  [\"69+28=\", \"32-8=\"],\n  [\"8+11=\", \"18+48=\"],\n  [\"78-65=\", \"28+54=\"],\n  [\"18-1=\", \"98-14=\"],\n  [\"68+25=\", \"88-52=\"],\n  [\"90-32=\", \"36+36=\"],\n  [\"49+16=\", \"17+0=\"],\n  [\"56-36=\", \"90-38=\"],\n  [\"78-15=\", \"16+14=\"],\n  [\"36-33=\", \"36+49=\"],\n  [\"86-20=\", \"1+91=\"],\n  [\"62-39=\", \"92-12=\"],\n  [\"12+84=\", \"34+48=\"],\n  [\"63-1=\", \"11+22=\"],\n  [\"69-42=\", \"79-76=\"],\n  [\"64-55=\", \"0+84=\"],\n  [\"81-66=\", \"19+46=\"],\n  [\"32+54=\", \"21+20=\"],\n  [\"92-87=\", \"80-45=\"],\n  [\"46-22=\", \"65-18=\"],\n  [\"50-26=\", \"38-19=\"],\n  [\"21+16=\", \"82-66=\"],\n  [\"26+29=\", \"84-20=\"],\n  [\"58-35=\", \"1+40=\"],\n  [\"45-22=\", \"42+9=\"],\n  [\"5+0=\", \"20+58=\"],\n  [\"79-61=\", \"21+19=\"],\n  [\"29+24=\", \"62+0=\"],\n  [\"97-38=\", \"93-85=\"],\n  [\"55+7=\", \"10+64=\"],\n  [\"69-19=\", \"77-41=\"],\n  [\"47+26=\", \"38+57=\"],\n  [\"79+15=\", \"45+21=\"],\n  [\"45-28=\", \"43+1=\"],\n  [\"20-2=\", \"85-52=\"],\n  [\"54-42=\", \"74-1=\"],\n  [\"18+50=\", \"67+24=\"],\n  [\"42-23=\", \"16-5=\"],\n  [\"28-13=\", \"2+20=\"],\n  [\"41-2=\", \"93-44=\"],\n  [\"70+26=\", \"34-10=\"],\n  [\"45-2=\", \"2+15=\"],\n  [\"96-2=\", \"78-22=\"],\n  [\"5+12=\", \"7+83=\"],\n  [\"75-24=\", \"66-16=\"],\n  [\"58-13=\", \"85-39=\"],\n  [\"81-28=\", \"93-34=\"],\n  [\"0+24=\", \"56-31=\"],\n  [\"29+68=\", \"89-66=\"],\n  [\"10+42=\", \"36+21=\"],\n  [\"0+21=\", \"89-7=\"],\n  [\"37+28=\", \"8+5=\"],\n  [\"16+11=\", \"2+70=\"],\n  [\"80-23=\", \"67+14=\"],\n  [\"51+20=\", \"22+5=\"],\n  [\"32+8=\", \"94-35=\"],\n  [\"28-11=\", \"9+44=\"],\n  [\"39-11=\", \"76-73=\"],\n  [\"35-32=\", \"87+3=\"],\n  [\"75-5=\", \"14+12=\"],\n  [\"87-82=\", \"38-27=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('37-31=', '83-9=')\n    ,@('42-22=', '80+16=')\n    ,@('42+35=', '74-23=')\n    ,@('99-55=', '91-61=')\n    ,@('8+19=', '16+20=')\n    ,@('70-32=', '11+3=')\n    ,@('52+15=', '12+23=')\n    ,@('16+46=', '28+1=')\n    ,@('4+80=', '24+62=')\n    ,@('52+18=', '31+58=')\n    ,@('42+41=', '56-49=')\n    ,@('11+86=', '71+23=')\n    ,@('23-14=', '67-14=')\n    ,@('69-33=', '53-25=')\n    ,@('56+39=', '99-11=')\n    ,@('27+13=', '20+34=')\n    ,@('68+6=', '8+67=')\n    ,@('74+19=', '84-83=')\n    ,@('44-11=', '47+7=')\n    ,@('62-21=', '73-45=')\n    ,@('90-69=', '47-32=')\n    ,@('6+60=', '81-78=')\n    ,@('58-29=', '81-74=')\n    ,@('76+11=', '4+9=')\n    ,@('81-15=', '19+26=')\n    ,@('90-6=', '13+53=')\n    ,@('8+44=', '87-37=')\n    ,@('42+12=', '29+4=')\n    ,@('77+3=', '36-19=')\n    ,@('77-35=', '58-54=')\n    ,@('85-44=', '83-19=')\n    ,@('66+28=', '89-2=')\n    ,@('37-12=', '73-41=')\n    ,@('3+3=', '52-50=')\n    ,@('13+76=', '40+59=')\n    ,@('56+23=', '54+2=')\n    ,@('27+71=', '21+0=')\n    ,@('53+25=', '8+51=')\n    ,@('64-44=', '61-50=')\n    ,@('69+28=', '32-8=')\n    ,@('8+11=', '18+48=')\n    ,@('78-65=', '28+54=')\n    ,@('18-1=', '98-14=')\n    ,@('68+25=', '88-52=')\n    ,@('90-32=', '36+36=')\n    ,@('49+16=', '17+0=')\n    ,@('56-36=', '90-38=')\n    ,@('78-15=', '16+14=')\n    ,@('36-33=', '36+49=')\n    ,@('86-20=', '1+91=')\n    ,@('62-39=', '92-12=')\n    ,@('12+84=', '34+48=')\n    ,@('63-1=', '11+22=')\n    ,@('69-42=', '79-76=')\n    ,@('64-55=', '0+84=')\n    ,@('81-66=', '19+46=')\n    ,@('32+54=', '21+20=')\n    ,@('92-87=', '80-45=')\n    ,@('46-22=', '65-18=')\n    ,@('50-26=', '38-19=')\n    ,@('21+16=', '82-66=')\n    ,@('26+29=', '84-20=')\n    ,@('58-35=', '1+40=')\n    ,@('45-22=', '42+9=')\n    ,@('5+0=', '20+58=')\n    ,@('79-61=', '21+19=')\n    ,@('29+24=', '62+0=')\n    ,@('97-38=', '93-85=')\n    ,@('55+7=', '10+64=')\n    ,@('69-19=', '77-41=')\n    ,@('47+26=', '38+57=')\n    ,@('79+15=', '45+21=')\n    ,@('45-28=', '43+1=')\n    ,@('20-2=', '85-52=')\n    ,@('54-42=', '74-1=')\n    ,@('18+50=', '67+24=')\n    ,@('42-23=', '16-5=')\n    ,@('28-13=', '2+20=')\n    ,@('41-2=', '93-44=')\n    ,@('70+26=', '34-10=')\n    ,@('45-2=', '2+15=')\n    ,@('96-2=', '78-22=')\n    ,@('5+12=', '7+83=')\n    ,@('75-24=', '66-16=')\n    ,@('58-13=', '85-39=')\n    ,@('81-28=', '93-34=')\n    ,@('0+24=', '56-31=')\n    ,@('29+68=', '89-66=')\n    ,@('10+42=', '36+21=')\n    ,@('0+21=', '89-7=')\n    ,@('37+28=', '8+5=')\n    ,@('16+11=', '2+70=')\n    ,@('80-23=', '67+14=')\n    ,@('51+20=', '22+5=')\n    ,@('32+8=', '94-35=')\n    ,@('28-11=', '9+44=')\n    ,@('39-11=', '76-73=')\n    ,@('35-32=', '87+3=')\n    ,@('75-5=', '14+12=')\n    ,@('87-82=', '38-27=')\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    [void]$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}"}
